# "Data input for trucks" — add an agnostic energy-consumption column to the
# storageAssets sheet so e-truck records (EV / EHGV) can carry a kWh/km figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("storageAssets")
$ws.Activate()

# New header in column O (shared string "energy_consumption_kwhpkm")
$ws.Range("O1").Value = "energy_consumption_kwhpkm"

# Per-row values: 0 for every existing storage asset, except the two
# e-vehicle rows which get their real-world consumption figures.
$values = @{
    2  = 0      # House_battery
    3  = 0      # House_heatmodel_A
    4  = 0      # House_heatmodel_B
    5  = 0      # House_heatmodel_C
    6  = 0      # House_heatmodel_D
    7  = 0      # House_heatmodel_E
    8  = 0      # House_heatmodel_F
    9  = 0      # House_heatmodel_G
    10 = 0.2    # EV
    11 = 0      # House_DH_heatdeliveryset
    12 = 0      # House_heatpump_MT_S
    13 = 0      # House_heatpump_MT_L
    14 = 0      # Grid_battery_1MWh
    15 = 1      # EHGV
}

foreach ($row in ($values.Keys | Sort-Object)) {
    $ws.Cells.Item($row, 15).Value = $values[$row]
}

# Match the author's final view state (scrolled right, last cell selected)
try {
    $excel.ActiveWindow.ScrollColumn = 5
} catch {
}
$ws.Range("O15").Select()
